$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WR")

# Add the new player's Week 17 row to the bottom of the WR stats table
$ws.Range("A10").Value = "K.Wilkerson"
$ws.Range("B10:J10").Value = 0

# Update the active selection to reflect where the user clicked next
$ws.Activate()
$ws.Range("K11").Select()
